{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document is a multiplication-practice worksheet: a date heading\n// paragraph followed by a 5-column table of \"A\u00d7B=C\" strings. This edit\n// updates the date and replaces each of the 100 table answer strings with\n// a new one (same shape, different numbers). Every \"old\" string below is\n// unique in the document, so an exact-text search-and-replace is safe and\n// unambiguous.\n\nconst PAIRS = [\n  [\"2023-04-27 Thursday\", \"2023-04-28 Friday\"],\n  [\"64\u00d732=2048\", \"16\u00d7100=1600\"],\n  [\"12\u00d739=468\", \"18\u00d737=666\"],\n  [\"67\u00d783=5561\", \"44\u00d762=2728\"],\n  [\"22\u00d765=1430\", \"33\u00d716=528\"],\n  [\"64\u00d727=1728\", \"42\u00d727=1134\"],\n  [\"77\u00d727=2079\", \"63\u00d754=3402\"],\n  [\"93\u00d772=6696\", \"27\u00d799=2673\"],\n  [\"66\u00d786=5676\", \"19\u00d724=456\"],\n  [\"26\u00d772=1872\", \"55\u00d745=2475\"],\n  [\"70\u00d797=6790\", \"26\u00d773=1898\"],\n  [\"29\u00d762=1798\", \"13\u00d765=845\"],\n  [\"30\u00d753=1590\", \"38\u00d754=2052\"],\n  [\"56\u00d728=1568\", \"84\u00d792=7728\"],\n  [\"19\u00d761=1159\", \"12\u00d786=1032\"],\n  [\"59\u00d728=1652\", \"42\u00d741=1722\"],\n  [\"92\u00d797=8924\", \"33\u00d794=3102\"],\n  [\"54\u00d770=3780\", \"20\u00d720=400\"],\n  [\"59\u00d775=4425\", \"76\u00d753=4028\"],\n  [\"38\u00d764=2432\", \"61\u00d710=610\"],\n  [\"88\u00d767=5896\", \"75\u00d766=4950\"],\n  [\"85\u00d715=1275\", \"66\u00d759=3894\"],\n  [\"27\u00d735=945\", \"55\u00d791=5005\"],\n  [\"43\u00d715=645\", \"88\u00d785=7480\"],\n  [\"50\u00d746=2300\", \"77\u00d762=4774\"],\n  [\"20\u00d718=360\", \"41\u00d747=1927\"],\n  [\"93\u00d794=8742\", \"84\u00d745=3780\"],\n  [\"59\u00d771=4189\", \"87\u00d799=8613\"],\n  [\"21\u00d798=2058\", \"80\u00d799=7920\"],\n  [\"60\u00d799=5940\", \"45\u00d737=1665\"],\n  [\"83\u00d744=3652\", \"73\u00d739=2847\"],\n  [\"38\u00d797=3686\", \"97\u00d782=7954\"],\n  [\"60\u00d793=5580\", \"57\u00d788=5016\"],\n  [\"91\u00d784=7644\", \"41\u00d795=3895\"],\n  [\"86\u00d732=2752\", \"95\u00d739=3705\"],\n  [\"47\u00d774=3478\", \"20\u00d778=1560\"],\n  [\"22\u00d787=1914\", \"27\u00d718=486\"],\n  [\"51\u00d735=1785\", \"40\u00d741=1640\"],\n  [\"33\u00d740=1320\", \"88\u00d712=1056\"],\n  [\"62\u00d725=1550\", \"80\u00d756=4480\"],\n  [\"38\u00d747=1786\", \"14\u00d727=378\"],\n  [\"77\u00d7100=7700\", \"67\u00d792=6164\"],\n  [\"68\u00d745=3060\", \"73\u00d782=5986\"],\n  [\"11\u00d731=341\", \"65\u00d743=2795\"],\n  [\"82\u00d732=2624\", \"31\u00d757=1767\"],\n  [\"83\u00d745=3735\", \"21\u00d731=651\"],\n  [\"81\u00d7100=8100\", \"38\u00d732=1216\"],\n  [\"61\u00d712=732\", \"70\u00d722=1540\"],\n  [\"13\u00d792=1196\", \"68\u00d774=5032\"],\n  [\"26\u00d733=858\", \"16\u00d782=1312\"],\n  [\"30\u00d731=930\", \"46\u00d736=1656\"],\n  [\"12\u00d716=192\", \"61\u00d714=854\"],\n  [\"79\u00d772=5688\", \"38\u00d779=3002\"],\n  [\"74\u00d789=6586\", \"25\u00d761=1525\"],\n  [\"66\u00d712=792\", \"24\u00d721=504\"],\n  [\"11\u00d765=715\", \"68\u00d714=952\"],\n  [\"90\u00d727=2430\", \"71\u00d760=4260\"],\n  [\"42\u00d774=3108\", \"48\u00d727=1296\"],\n  [\"50\u00d723=1150\", \"65\u00d758=3770\"],\n  [\"92\u00d764=5888\", \"62\u00d722=1364\"],\n  [\"34\u00d737=1258\", \"57\u00d746=2622\"],\n  [\"68\u00d710=680\", \"40\u00d7100=4000\"],\n  [\"33\u00d719=627\", \"91\u00d781=7371\"],\n  [\"16\u00d793=1488\", \"95\u00d748=4560\"],\n  [\"91\u00d762=5642\", \"11\u00d746=506\"],\n  [\"22\u00d777=1694\", \"38\u00d711=418\"],\n  [\"77\u00d756=4312\", \"30\u00d773=2190\"],\n  [\"47\u00d720=940\", \"59\u00d724=1416\"],\n  [\"17\u00d775=1275\", \"57\u00d756=3192\"],\n  [\"77\u00d797=7469\", \"29\u00d777=2233\"],\n  [\"34\u00d779=2686\", \"66\u00d783=5478\"],\n  [\"47\u00d768=3196\", \"38\u00d760=2280\"],\n  [\"95\u00d781=7695\", \"84\u00d754=4536\"],\n  [\"69\u00d716=1104\", \"31\u00d743=1333\"],\n  [\"31\u00d720=620\", \"49\u00d788=4312\"],\n  [\"90\u00d748=4320\", \"70\u00d715=1050\"],\n  [\"19\u00d799=1881\", \"57\u00d759=3363\"],\n  [\"69\u00d766=4554\", \"51\u00d783=4233\"],\n  [\"51\u00d727=1377\", \"79\u00d742=3318\"],\n  [\"17\u00d732=544\", \"68\u00d7100=6800\"],\n  [\"51\u00d796=4896\", \"96\u00d796=9216\"],\n  [\"44\u00d783=3652\", \"20\u00d717=340\"],\n  [\"66\u00d728=1848\", \"54\u00d743=2322\"],\n  [\"97\u00d770=6790\", \"56\u00d729=1624\"],\n  [\"56\u00d735=1960\", \"23\u00d770=1610\"],\n  [\"23\u00d762=1426\", \"42\u00d716=672\"],\n  [\"93\u00d728=2604\", \"97\u00d741=3977\"],\n  [\"97\u00d723=2231\", \"26\u00d716=416\"],\n  [\"67\u00d728=1876\", \"39\u00d757=2223\"],\n  [\"13\u00d760=780\", \"93\u00d722=2046\"],\n  [\"16\u00d753=848\", \"39\u00d735=1365\"],\n  [\"12\u00d779=948\", \"31\u00d789=2759\"],\n  [\"41\u00d783=3403\", \"17\u00d750=850\"],\n  [\"39\u00d724=936\", \"98\u00d789=8722\"],\n  [\"80\u00d779=6320\", \"89\u00d729=2581\"],\n  [\"51\u00d755=2805\", \"85\u00d792=7820\"],\n  [\"63\u00d765=4095\", \"10\u00d743=430\"],\n  [\"17\u00d762=1054\", \"70\u00d789=6230\"],\n  [\"93\u00d759=5487\", \"21\u00d757=1197\"],\n  [\"92\u00d722=2024\", \"51\u00d739=1989\"],\n  [\"50\u00d741=2050\", \"21\u00d788=1848\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of PAIRS) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found, expected exactly one match: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The document is a multiplication-practice worksheet: a date heading\n# paragraph followed by a 5-column table of \"A x B = C\" answer strings.\n# This script updates the date and replaces each of the 100 table answer\n# strings with a new one (same shape, different numbers). Every \"old\"\n# string is unique in the document, so an exact Find/Replace (whole\n# document, match case) is safe and unambiguous for each pair.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2023-04-27 Thursday\", \"2023-04-28 Friday\")\n    ,@(\"64\u00d732=2048\", \"16\u00d7100=1600\")\n    ,@(\"12\u00d739=468\", \"18\u00d737=666\")\n    ,@(\"67\u00d783=5561\", \"44\u00d762=2728\")\n    ,@(\"22\u00d765=1430\", \"33\u00d716=528\")\n    ,@(\"64\u00d727=1728\", \"42\u00d727=1134\")\n    ,@(\"77\u00d727=2079\", \"63\u00d754=3402\")\n    ,@(\"93\u00d772=6696\", \"27\u00d799=2673\")\n    ,@(\"66\u00d786=5676\", \"19\u00d724=456\")\n    ,@(\"26\u00d772=1872\", \"55\u00d745=2475\")\n    ,@(\"70\u00d797=6790\", \"26\u00d773=1898\")\n    ,@(\"29\u00d762=1798\", \"13\u00d765=845\")\n    ,@(\"30\u00d753=1590\", \"38\u00d754=2052\")\n    ,@(\"56\u00d728=1568\", \"84\u00d792=7728\")\n    ,@(\"19\u00d761=1159\", \"12\u00d786=1032\")\n    ,@(\"59\u00d728=1652\", \"42\u00d741=1722\")\n    ,@(\"92\u00d797=8924\", \"33\u00d794=3102\")\n    ,@(\"54\u00d770=3780\", \"20\u00d720=400\")\n    ,@(\"59\u00d775=4425\", \"76\u00d753=4028\")\n    ,@(\"38\u00d764=2432\", \"61\u00d710=610\")\n    ,@(\"88\u00d767=5896\", \"75\u00d766=4950\")\n    ,@(\"85\u00d715=1275\", \"66\u00d759=3894\")\n    ,@(\"27\u00d735=945\", \"55\u00d791=5005\")\n    ,@(\"43\u00d715=645\", \"88\u00d785=7480\")\n    ,@(\"50\u00d746=2300\", \"77\u00d762=4774\")\n    ,@(\"20\u00d718=360\", \"41\u00d747=1927\")\n    ,@(\"93\u00d794=8742\", \"84\u00d745=3780\")\n    ,@(\"59\u00d771=4189\", \"87\u00d799=8613\")\n    ,@(\"21\u00d798=2058\", \"80\u00d799=7920\")\n    ,@(\"60\u00d799=5940\", \"45\u00d737=1665\")\n    ,@(\"83\u00d744=3652\", \"73\u00d739=2847\")\n    ,@(\"38\u00d797=3686\", \"97\u00d782=7954\")\n    ,@(\"60\u00d793=5580\", \"57\u00d788=5016\")\n    ,@(\"91\u00d784=7644\", \"41\u00d795=3895\")\n    ,@(\"86\u00d732=2752\", \"95\u00d739=3705\")\n    ,@(\"47\u00d774=3478\", \"20\u00d778=1560\")\n    ,@(\"22\u00d787=1914\", \"27\u00d718=486\")\n    ,@(\"51\u00d735=1785\", \"40\u00d741=1640\")\n    ,@(\"33\u00d740=1320\", \"88\u00d712=1056\")\n    ,@(\"62\u00d725=1550\", \"80\u00d756=4480\")\n    ,@(\"38\u00d747=1786\", \"14\u00d727=378\")\n    ,@(\"77\u00d7100=7700\", \"67\u00d792=6164\")\n    ,@(\"68\u00d745=3060\", \"73\u00d782=5986\")\n    ,@(\"11\u00d731=341\", \"65\u00d743=2795\")\n    ,@(\"82\u00d732=2624\", \"31\u00d757=1767\")\n    ,@(\"83\u00d745=3735\", \"21\u00d731=651\")\n    ,@(\"81\u00d7100=8100\", \"38\u00d732=1216\")\n    ,@(\"61\u00d712=732\", \"70\u00d722=1540\")\n    ,@(\"13\u00d792=1196\", \"68\u00d774=5032\")\n    ,@(\"26\u00d733=858\", \"16\u00d782=1312\")\n    ,@(\"30\u00d731=930\", \"46\u00d736=1656\")\n    ,@(\"12\u00d716=192\", \"61\u00d714=854\")\n    ,@(\"79\u00d772=5688\", \"38\u00d779=3002\")\n    ,@(\"74\u00d789=6586\", \"25\u00d761=1525\")\n    ,@(\"66\u00d712=792\", \"24\u00d721=504\")\n    ,@(\"11\u00d765=715\", \"68\u00d714=952\")\n    ,@(\"90\u00d727=2430\", \"71\u00d760=4260\")\n    ,@(\"42\u00d774=3108\", \"48\u00d727=1296\")\n    ,@(\"50\u00d723=1150\", \"65\u00d758=3770\")\n    ,@(\"92\u00d764=5888\", \"62\u00d722=1364\")\n    ,@(\"34\u00d737=1258\", \"57\u00d746=2622\")\n    ,@(\"68\u00d710=680\", \"40\u00d7100=4000\")\n    ,@(\"33\u00d719=627\", \"91\u00d781=7371\")\n    ,@(\"16\u00d793=1488\", \"95\u00d748=4560\")\n    ,@(\"91\u00d762=5642\", \"11\u00d746=506\")\n    ,@(\"22\u00d777=1694\", \"38\u00d711=418\")\n    ,@(\"77\u00d756=4312\", \"30\u00d773=2190\")\n    ,@(\"47\u00d720=940\", \"59\u00d724=1416\")\n    ,@(\"17\u00d775=1275\", \"57\u00d756=3192\")\n    ,@(\"77\u00d797=7469\", \"29\u00d777=2233\")\n    ,@(\"34\u00d779=2686\", \"66\u00d783=5478\")\n    ,@(\"47\u00d768=3196\", \"38\u00d760=2280\")\n    ,@(\"95\u00d781=7695\", \"84\u00d754=4536\")\n    ,@(\"69\u00d716=1104\", \"31\u00d743=1333\")\n    ,@(\"31\u00d720=620\", \"49\u00d788=4312\")\n    ,@(\"90\u00d748=4320\", \"70\u00d715=1050\")\n    ,@(\"19\u00d799=1881\", \"57\u00d759=3363\")\n    ,@(\"69\u00d766=4554\", \"51\u00d783=4233\")\n    ,@(\"51\u00d727=1377\", \"79\u00d742=3318\")\n    ,@(\"17\u00d732=544\", \"68\u00d7100=6800\")\n    ,@(\"51\u00d796=4896\", \"96\u00d796=9216\")\n    ,@(\"44\u00d783=3652\", \"20\u00d717=340\")\n    ,@(\"66\u00d728=1848\", \"54\u00d743=2322\")\n    ,@(\"97\u00d770=6790\", \"56\u00d729=1624\")\n    ,@(\"56\u00d735=1960\", \"23\u00d770=1610\")\n    ,@(\"23\u00d762=1426\", \"42\u00d716=672\")\n    ,@(\"93\u00d728=2604\", \"97\u00d741=3977\")\n    ,@(\"97\u00d723=2231\", \"26\u00d716=416\")\n    ,@(\"67\u00d728=1876\", \"39\u00d757=2223\")\n    ,@(\"13\u00d760=780\", \"93\u00d722=2046\")\n    ,@(\"16\u00d753=848\", \"39\u00d735=1365\")\n    ,@(\"12\u00d779=948\", \"31\u00d789=2759\")\n    ,@(\"41\u00d783=3403\", \"17\u00d750=850\")\n    ,@(\"39\u00d724=936\", \"98\u00d789=8722\")\n    ,@(\"80\u00d779=6320\", \"89\u00d729=2581\")\n    ,@(\"51\u00d755=2805\", \"85\u00d792=7820\")\n    ,@(\"63\u00d765=4095\", \"10\u00d743=430\")\n    ,@(\"17\u00d762=1054\", \"70\u00d789=6230\")\n    ,@(\"93\u00d759=5487\", \"21\u00d757=1197\")\n    ,@(\"92\u00d722=2024\", \"51\u00d739=1989\")\n    ,@(\"50\u00d741=2050\", \"21\u00d788=1848\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $oldText,  # FindText\n        $true,     # MatchCase\n        $true,     # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap -> wdFindContinue\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace -> wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Text not found, expected exactly one match: $oldText\"\n    }\n}\n"}
